# Apply updates to column F ("dSF") values on specific rows of Sheet1.
# Mapping of row -> new value, as derived from the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    7  = 1
    8  = 1
    9  = -1
    15 = 1
    19 = 2
    20 = 3
    24 = 0
    27 = 2
    35 = 0
    36 = 1
    42 = 0
    43 = 3
    48 = 0
    57 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

$wb.Save()
